# Include the image title in the picture's description (AlternativeText).
# Pandoc writes the markdown image link into the pptx <p:cNvPr descr="...">
# attribute; when the image also carries a title (as in
# `![alt text](link "title")`), the title should be included too. Pandoc
# signals a "title" image with the literal "fig:" prefix. Every picture in
# this deck was inserted from a titled image link, so update each picture's
# AlternativeText (which maps to descr) from "lalune.jpg" to
# "fig:  lalune.jpg".

$p = $ppt.ActivePresentation

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    for ($j = 1; $j -le $s.Shapes.Count; $j++) {
        $shp = $s.Shapes.Item($j)
        if ($shp.Type -eq 13 -and $shp.AlternativeText -eq "lalune.jpg") {
            $shp.AlternativeText = "fig:  lalune.jpg"
        }
    }
}
